$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.059.57'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '3.471.79'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  +0.10%  '
$style_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.27'
$ws.Range('D5').Style = $style_D5
$ws.Range('E5').Value = '  -0.50%  '
$style_D6 = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.72'
$ws.Range('D6').Style = $style_D6
$ws.Range('E6').Value = '  -2.65%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -1.10%  '
$style_D9 = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.67'
$ws.Range('D9').Style = $style_D9
$ws.Range('E9').Value = '  +5.67%  '
$ws.Range('E10').Value = '  -0.96%  '
$style_D11 = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.384'
$ws.Range('D11').Style = $style_D11
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '4.068.94'
$ws.Range('E12').Value = '  -0.33%  '
$style_D13 = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.119'
$ws.Range('D13').Style = $style_D13
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('E14').Value = '  -2.46%  '
$ws.Range('D15').Value = '3.468.09'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').Value = '64.061.67'
$ws.Range('E16').Value = '  -0.45%  '
$style_D17 = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.93'
$ws.Range('D17').Style = $style_D17
$ws.Range('E17').Value = '  -3.16%  '
$style_D18 = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.94'
$ws.Range('D18').Style = $style_D18
$ws.Range('E18').Value = '  +0.64%  '
$style_D19 = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.67'
$ws.Range('D19').Style = $style_D19
$ws.Range('E19').Value = '  -1.53%  '
$style_D20 = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.36'
$ws.Range('D20').Style = $style_D20
$ws.Range('E20').Value = '  -1.90%  '
$style_D21 = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '384.63'
$ws.Range('D21').Style = $style_D21
$ws.Range('E21').Value = '  -2.46%  '
$style_D22 = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.565'
$ws.Range('D22').Style = $style_D22
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').Value = '3.614.39'
$ws.Range('E23').Value = '  -0.45%  '
$style_D24 = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.44'
$ws.Range('D24').Style = $style_D24
$ws.Range('E24').Value = '  -0.46%  '
$style_D25 = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('D25').Style = $style_D25
$ws.Range('E25').Value = '  -0.24%  '
$style_D26 = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.39'
$ws.Range('D26').Style = $style_D26
$ws.Range('E26').Value = '  -6.04%  '
$style_D27 = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000111'
$ws.Range('D27').Style = $style_D27
$ws.Range('E27').Value = '  -3.16%  '
$style_D28 = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.999'
$ws.Range('D28').Style = $style_D28
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('E29').Value = '  -0.53%  '
$style_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.06'
$ws.Range('D30').Style = $style_D30
$ws.Range('E30').Value = '  -4.06%  '
$ws.Range('E31').Value = '  +3.69%  '
$style_D32 = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.93'
$ws.Range('D32').Style = $style_D32
$ws.Range('E32').Value = '  -3.55%  '
$style_D33 = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.43'
$ws.Range('D33').Style = $style_D33
$ws.Range('E33').Value = '  -3.97%  '
$ws.Range('D34').Value = '3.502.55'
$ws.Range('E34').Value = '  -0.32%  '
$ws.Range('E35').Value = '  -0.07%  '
$style_D36 = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.91'
$ws.Range('D36').Style = $style_D36
$ws.Range('E36').Value = '  -2.08%  '
$style_D37 = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.18'
$ws.Range('D37').Style = $style_D37
$ws.Range('E37').Value = '  +1.11%  '
$style_D38 = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.75'
$ws.Range('D38').Style = $style_D38
$ws.Range('E38').Value = '  -2.03%  '
$style_D39 = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '163.00'
$ws.Range('D39').Style = $style_D39
$ws.Range('E39').Value = '  -1.92%  '
$style_D40 = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.49'
$ws.Range('D40').Style = $style_D40
$ws.Range('E40').Value = '  -3.76%  '
$ws.Range('E41').Value = '  -0.87%  '
$style_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.795'
$ws.Range('D42').Style = $style_D42
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('E43').Value = '  +0.13%  '
$style_D44 = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.51'
$ws.Range('D44').Style = $style_D44
$ws.Range('E44').Value = '  -0.77%  '
$style_D45 = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.34'
$ws.Range('D45').Style = $style_D45
$ws.Range('E45').Value = '  -1.04%  '
$style_D46 = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.61'
$ws.Range('D46').Style = $style_D46
$ws.Range('E46').Value = '  -2.08%  '
$style_D47 = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.46'
$ws.Range('D47').Style = $style_D47
$ws.Range('E47').Value = '  -6.75%  '
$ws.Range('E48').Value = '  -3.92%  '
$style_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.70'
$ws.Range('D49').Style = $style_D49
$ws.Range('E49').Value = '  -0.86%  '
$style_D50 = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.893'
$ws.Range('D50').Style = $style_D50
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = '2.330.79'
$ws.Range('E51').Value = '  -5.08%  '
